$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (+ B2) first, matching the order strings were authored in the
# workbook: title, then Name/Value header, then the command-name list.
$ws.Range("A1").Value = "xArm App MainWindow.cs"
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "Value"

$ws.Range("A3").Value = "CMD_MULT_SERVO_MOVE"
$ws.Range("B3").Value = 3

$ws.Range("A4").Value = "CMD_ACTION_DOWNLOAD"
$ws.Range("B4").Value = 5

$ws.Range("A5").Value = "CMD_FULL_ACTION_RUN"
$ws.Range("B5").Value = 6

$ws.Range("A6").Value = "CMD_FULL_ACTION_STOP"
$ws.Range("B6").Value = 7

$ws.Range("A7").Value = "CMD_FULL_ACTION_ERASE"
$ws.Range("B7").Value = 8

$ws.Range("A8").Value = "CMD_SERVO_OFFSET_WRITE"
$ws.Range("B8").Value = 12

$ws.Range("A9").Value = "CMD_SERVO_OFFSET_READ"
$ws.Range("B9").Value = 13

$ws.Range("A10").Value = "CMD_SERVO_OFFSET_ADJUST"
$ws.Range("B10").Value = 14

$ws.Range("A11").Value = "CMD_MULT_SERVO_UNLOAD"
$ws.Range("B11").Value = 20

$ws.Range("A12").Value = "CMD_MULT_SERVO_POS_READ"
$ws.Range("B12").Value = 21

$ws.Range("A13").Value = "CMD_BUS_SERVO_OFFSET_WRITE"
$ws.Range("B13").Value = 22

$ws.Range("A14").Value = "CMD_BUS_SERVO_OFFSET_READ"
$ws.Range("B14").Value = 23

$ws.Range("A15").Value = "CMD_BUS_SERVO_OFFSET_ADJUST"
$ws.Range("B15").Value = 24

$ws.Range("A16").Value = "CMD_BUS_SERVO_MOROR_CTRL"
$ws.Range("B16").Value = 26

$ws.Range("A17").Value = "CMD_BUS_SERVO_INFO_WRITE"
$ws.Range("B17").Value = 27

$ws.Range("A18").Value = "CMD_BUS_SERVO_INFO_READ"
$ws.Range("B18").Value = 28

# Small "Command Format" diagram off to the right.
$ws.Range("C2").Value = "Length"

$ws.Range("G2").Value = "Command Format"

$ws.Range("G3").Value = "Header"
$ws.Range("H3").Value = "Data Length"
$ws.Range("I3").Value = "Command"
$ws.Range("J3").Value = "Parameter"

$ws.Range("G4").Value = "0x55 0x55"
$ws.Range("H4").Value = "Length"
$ws.Range("I4").Value = "Cmd"
$ws.Range("J4").Value = "Prm 1…Prm N"

$ws.Columns.Item(1).ColumnWidth = 32.140625
$ws.Columns.Item(7).ColumnWidth = 10.7109375
$ws.Columns.Item(8).ColumnWidth = 12.42578125
$ws.Columns.Item(9).ColumnWidth = 11
$ws.Columns.Item(10).ColumnWidth = 15.140625

$ws.Range("J8").Select()
